$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the team-members block
#    (it sat right after "A01007088").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Insert the new "Instructions" block just before the final blank
#    paragraph that precedes the section break.
# ------------------------------------------------------------------
$pLast = $d.Paragraphs.Last
$null = $pLast.Range.InsertParagraphBefore()
$null = $pLast.Range.InsertParagraphBefore()
$null = $pLast.Range.InsertParagraphBefore()

$pLastIndex = $d.Paragraphs.Count
$pInstructions = $d.Paragraphs.Item($pLastIndex - 3)
$pCompile      = $d.Paragraphs.Item($pLastIndex - 2)
$pSassCmd      = $d.Paragraphs.Item($pLastIndex - 1)

$pInstructions.Range.Text = "Instructions:"
$pInstructions.Range.Bold = 1
$pInstructions.Range.Font.Underline = 1

$pCompile.Range.Text = "To compile Sass:"
$pCompile.Range.Bold = 1

$sassXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">sass --watch --style expanded </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>scss:styles</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$pSassCmd.Range.InsertXML($sassXml)

# ------------------------------------------------------------------
# 3. Re-add a trailing paragraph holding the (now relocated) "_GoBack"
#    bookmark as the very last paragraph of the document body.
# ------------------------------------------------------------------
$pEnd = $d.Paragraphs.Last
$null = $pEnd.Range.InsertParagraphAfter()
$pBookmark = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $pBookmark.Range)
